# Fixed some bugs in gamecollection: correct shuffled row data in Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (A, B, C, D, E, F) target values
$rows = @{
    3  = @(701, 3, 90, 45, 97, 15)
    4  = @(1202, 2, 10, 10, 10, 10)
    5  = @(101, 9, 30, 15, 60, 15)
    7  = @(401, 9, 48, 67, 75, 45)
    8  = @(801, 3, 67, 65, 52, 45)
    9  = @(201, 9, 30, 15, 45, 30)
    10 = @(1001, 18, 30, 75, 60, 72)
    11 = @(301, 6, 45, 30, 60, 45)
    12 = @(501, 9, 52, 30, 75, 45)
    13 = @(601, 9, 60, 67, 60, 42)
    14 = @(1203, 3, 15, 15, 15, 15)
    15 = @(901, 16, 15, 45, 60, 60)
    16 = @(2, 0, 2, 2, 2, 2)
    17 = @(3, 0, 3, 3, 3, 3)
    18 = @(502, 0, 4, 0, 0, 0)
    20 = @(1, 0, 2, 2, 2, 2)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
